$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent results for case "380 kV" (rows 2-25, cols B-F and J-O)
# Each entry: row number, then B,C,D,E,F values, then J,K,L,M,N,O values
$data = @(
    @(2, 10.72688076740759,3.963739053596491,8.828761819488594,13.59623682752078,34.8337867893789, 10.19193600088265,10.14184448770887,0,15.3995089895444,20.24264389651791,26.46393044682838),
    @(3, 10.50258162370882,3.79030017283605,8.799366737831663,13.59191496506123,34.88930835545293, 10.2112871807631,9.996660893980611,0,15.3467258366097,20.30316119320692,26.53363303243357),
    @(4, 10.36454253975447,3.678900705130868,8.782682844002618,13.59151615728081,34.93059362212934, 10.22429469031937,9.908205394185815,0,15.31668049222607,20.34205809608637,26.58128892995628),
    @(5, 10.30829268747907,3.6323079961131,8.776232209314111,13.59192231188345,34.949225018105, 10.22987876008328,9.872375087786041,0,15.30504012954249,20.35834747052239,26.60192932555659),
    @(6, 10.29895501941307,3.62450025456351,8.775182255342415,13.59202414393033,34.95242784061854, 10.23082311650926,9.866439796619368,0,15.30314394833782,20.36107883485682,26.60543030325809),
    @(7, 10.36378381545839,3.678277130484823,8.782594432334847,13.59151933026123,34.93083757673993, 10.22436885113644,9.907721242636516,0,15.31652105182465,20.34227600284247,26.58156235483321),
    @(8, 10.64965721480564,3.904976560875141,8.818346226665497,13.594279736216,34.8514363485429, 10.19837481503342,10.09166828091744,0,15.38082432969454,20.26315009451844,26.4869547542282),
    @(9, 11.20423407606156,4.309226284384185,8.899050801489091,13.6175010819387,34.75288567505802, 10.15631957874148,10.455957063392,0,15.52525845548938,20.12172366107787,26.34003751757073),
    @(10, 11.60354828633749,4.580156980640028,8.96447079527619,13.64529332349841,34.71539096037279, 10.13084010677883,10.72327486676924,0,15.64195538442661,20.02610902314032,26.2557115554231),
    @(11, 11.78256805715212,4.697513649901563,8.995485210240615,13.66023501527034,34.70591804519916, 10.1204213131959,10.84429882644513,0,15.69720354020913,19.98439380994133,26.22249148375045),
    @(12, 11.84991279533181,4.741091074227562,9.007403007657302,13.6662204374827,34.70342071162595, 10.11664416910625,10.89000250858788,0,15.71842307974655,19.96885202963475,26.210652042093),
    @(13, 11.83542977930713,4.731744501398332,9.004828704136758,13.66491686625289,34.70391010226236, 10.11745016756501,10.88016563671549,0,15.71384001287363,19.97218791306409,26.21316893728546),
    @(14, 11.78811786461989,4.701116167576213,8.996462262657227,13.66072089896251,34.70569075162065, 10.12010719549654,10.84806164597024,0,15.69894338084717,19.98311007697813,26.22150260356739),
    @(15, 11.75907790561414,4.682242617759473,8.991359931379534,13.65819327468698,34.70692335237804, 10.12175660095209,10.82837946268008,0,15.68985722341057,19.98983337485806,26.22670365172017),
    @(16, 11.59179008875745,4.572367480890123,8.962468574458061,13.64436283645858,34.71616259613874, 10.13154455676858,10.7153503655868,0,15.63838722367265,20.02887093748389,26.25798609873907),
    @(17, 11.48844104139117,4.503441585041039,8.945061015126694,13.63646488031498,34.72377263248593, 10.13784910530718,10.64583171712846,0,15.60735695874916,20.05327439122058,26.27849424431241),
    @(18, 11.42875362001473,4.463243489770009,8.935167170673532,13.63213877595424,34.72886360613163, 10.14158564233409,10.60579424414634,0,15.58971346744294,20.06747825672864,26.29077375187997),
    @(19, 11.40850475611762,4.449538530390309,8.931837853289174,13.63071132557144,34.73070994295995, 10.1428697273335,10.59223058435555,0,15.58377515089673,20.07231626980418,26.29501443281083),
    @(20, 11.49946846048265,4.510836298617345,8.946901863384788,13.63728323752571,34.72288864995502, 10.13716655900162,10.65323780581659,0,15.61063912979161,20.05065925969363,26.27626104081032),
    @(21, 11.80202714924935,4.71013598581712,8.998915043406388,13.66194449987866,34.70513816156951, 10.11932219877144,10.85749509910929,0,15.70331089094392,19.97989506480246,26.21903470460352),
    @(22, 11.99713775677318,4.83535416869263,9.033915399823227,13.6799682838295,34.69988897490124, 10.10864028853763,10.99023930195388,0,15.76560958579149,19.93513159891778,26.1859493137072),
    @(23, 11.89326541309087,4.768988169714834,9.015145338402807,13.67017537276748,34.70210974489471, 10.11425181806766,10.91947332780525,0,15.73220537699982,19.95888719990926,26.20321242371322),
    @(24, 11.49448380219082,4.507494928053963,8.946069260285974,13.63691258974044,34.72328606854605, 10.13747478939421,10.64988972879181,0,15.60915464809843,20.05184101857312,26.27726914866744),
    @(25, 11.05532522304545,4.204352550518693,8.876116895507717,13.609324425221,34.77341839042774, 10.16674370238564,10.35728600780902,0,15.48428168940532,20.15852114709146,26.37564078331037)
)

foreach ($row in $data) {
    $r = $row[0]

    $bf = New-Object 'object[,]' 1,5
    for ($i = 0; $i -lt 5; $i++) { $bf[0,$i] = $row[1 + $i] }
    $ws.Range("B${r}:F${r}").Value = $bf

    $jo = New-Object 'object[,]' 1,6
    for ($i = 0; $i -lt 6; $i++) { $jo[0,$i] = $row[6 + $i] }
    $ws.Range("J${r}:O${r}").Value = $jo
}

Write-Output "Updated loading_percent values for rows 2-25"